# Simulation_Parameters.xlsx update
# - Adds a "Value" column header to the Initialization_Parameters and
#   Main_Loop_Parameters sheets.
# - Populates the previously-empty Randomness_Parameters sheet with the new
#   "Agent Randomness Parameters" table (bank/household/consumer firm/
#   capital firm randomness seeds) and makes it the active/selected sheet.
# - Updates the remembered cell selection on each sheet to match where the
#   cursor was left after the edits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Randomness_Parameters
# (populated first so the new shared strings land in the same order as
# the target workbook: "Agent Randomness Parameters", bank_randomness,
# household_randomness, consumer_firm_randomness, capital_firm_randomness,
# then "Value" last, reused afterwards by the other two sheets)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Randomness_Parameters")

$ws3.Range("A1").Value = "Agent Randomness Parameters"
$ws3.Range("A2").Value = "Bank"
$ws3.Range("A3").Value = "bank_randomness"
$ws3.Range("B3").Value = 1
$ws3.Range("A4").Value = "Households"
$ws3.Range("A5").Value = "household_randomness"
$ws3.Range("B5").Value = 2
$ws3.Range("A6").Value = "Consumer Firms"
$ws3.Range("A7").Value = "consumer_firm_randomness"
$ws3.Range("B7").Value = 3
$ws3.Range("A8").Value = "Capital Firms"
$ws3.Range("A9").Value = "capital_firm_randomness"
$ws3.Range("B9").Value = 0
$ws3.Range("B1").Value = "Value"

$ws3.Range("A1:B1").Font.Bold = $true
$ws3.Range("A2").Font.Bold = $true
$ws3.Range("A4").Font.Bold = $true
$ws3.Range("A6").Font.Bold = $true
$ws3.Range("A8").Font.Bold = $true

$ws3.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Sheet: Initialization_Parameters
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Initialization_Parameters")
$ws1.Range("B1").Value = "Value"
$ws1.Range("B1").Font.Bold = $true
$ws1.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: Main_Loop_Parameters
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Main_Loop_Parameters")
$ws2.Range("B1").Value = "Value"
$ws2.Range("B1").Font.Bold = $true
$ws2.Range("D26").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: Consumer_Sectors
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Consumer_Sectors")
$ws4.Range("I11").Select() | Out-Null

# ---------------------------------------------------------------------
# Final selection / active sheet: Randomness_Parameters, cell G12
# ---------------------------------------------------------------------
$ws3.Range("G12").Select() | Out-Null
$ws3.Activate()
